$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "Sprite/well"
$ws.Range("I2").Value = "Prefabs/well"
$ws.Range("H3").Value = "Sprite/well"
$ws.Range("I3").Value = "Prefabs/well"
$ws.Range("H4").Value = "Sprite/well"
$ws.Range("I4").Value = "Prefabs/well"

$ws.Range("J3").Select()
